# تعديل حدث في Card24 - الصف 23
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Row 21: fill previously-empty cells with "nan"
$ws.Range("A21").Value = "nan"
$ws.Range("L21").Value = "nan"
$ws.Range("M21").Value = "nan"
$ws.Range("N21").Value = "nan"

# Row 25: clear the event data (shifted up to row 21 conceptually), leaving blanks
$ws.Range("A25").Value = ""
$ws.Range("L25").Value = ""
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = ""
